$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览): refresh "want to go" counts (column F) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 883
$ws1.Range("F3").Value = 1464
$ws1.Range("F4").Value = 1121
$ws1.Range("F8").Value = 684
$ws1.Range("F11").Value = 98
$ws1.Range("F12").Value = 223
$ws1.Range("F13").Value = 161
$ws1.Range("F14").Value = 2883
$ws1.Range("F15").Value = 13
$ws1.Range("F16").Value = 10
$ws1.Range("F17").Value = 438
$ws1.Range("F19").Value = 509
$ws1.Range("F20").Value = 284
$ws1.Range("F24").Value = 676
$ws1.Range("F25").Value = 58
$ws1.Range("F26").Value = 253
$ws1.Range("F27").Value = 970
$ws1.Range("F29").Value = 1598
$ws1.Range("F30").Value = 329

# --- Sheet 2 (演出): refresh "want to go" counts (column F) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 33
$ws2.Range("F4").Value = 656
$ws2.Range("F5").Value = 237
$ws2.Range("F6").Value = 27
$ws2.Range("F7").Value = 233
$ws2.Range("F9").Value = 67
$ws2.Range("F10").Value = 43
$ws2.Range("F12").Value = 133

# --- Sheet 3 (本地生活): refresh "want to go" counts (column F) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 68

# --- Sheet 4 (全部类型): refresh "want to go" counts for rows 3-37 (mirrors sheets 1-3) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 883
$ws4.Range("F4").Value = 1464
$ws4.Range("F5").Value = 1121
$ws4.Range("F7").Value = 33
$ws4.Range("F8").Value = 68
$ws4.Range("F12").Value = 684
$ws4.Range("F13").Value = 656
$ws4.Range("F16").Value = 98
$ws4.Range("F17").Value = 223
$ws4.Range("F18").Value = 161
$ws4.Range("F19").Value = 2883
$ws4.Range("F20").Value = 13
$ws4.Range("F21").Value = 10
$ws4.Range("F22").Value = 237
$ws4.Range("F23").Value = 438
$ws4.Range("F25").Value = 509
$ws4.Range("F26").Value = 284
$ws4.Range("F28").Value = 27
$ws4.Range("F31").Value = 233
$ws4.Range("F33").Value = 67
$ws4.Range("F34").Value = 676
$ws4.Range("F35").Value = 43
$ws4.Range("F37").Value = 133

# --- Sheet 4: remove duplicate "angela LIVE" row (old row 38) by shifting rows 39-45 up into 38-44 ---
# (content only - B:I - the "A" index column is left untouched; row 45 is removed afterwards)
$ws4.Range("B38").Value = '2024-05-04'
$ws4.Range("C38").Value = '广州·运动番ONLY'
$ws4.Range("D38").Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws4.Range("E38").Value = '2024.05.04 10:00-05.04 17:00'
$ws4.Range("F38").Value = 58
$ws4.Range("G38").Value = 60
$ws4.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=82526'
$ws4.Range("I38").Value = '//i1.hdslb.com/bfs/openplatform/202403/CawIgD2O1709803813638.jpeg'

$ws4.Range("B39").Value = '2024-05-04'
$ws4.Range("C39").Value = '广州·黑塔利亚Only'
$ws4.Range("D39").Value = '迎宾大道123号 赛仑吉地大酒店'
$ws4.Range("E39").Value = '2024.05.04 09:30-05.04 16:00'
$ws4.Range("F39").Value = 253
$ws4.Range("G39").Value = 68
$ws4.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=82056'
$ws4.Range("I39").Value = '//i2.hdslb.com/bfs/openplatform/202402/KI6tnMd81708917202487.jpeg'

$ws4.Range("B40").Value = '2024-05-05'
$ws4.Range("C40").Value = '广州·第八届萌物语动漫嘉年华'
$ws4.Range("D40").Value = '洛浦街厦滘西环路1号 岭南会展中心'
$ws4.Range("E40").Value = '2024.05.05 10:00-05.05 17:00'
$ws4.Range("F40").Value = 970
$ws4.Range("G40").Value = 60
$ws4.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=81566'
$ws4.Range("I40").Value = '//i2.hdslb.com/bfs/openplatform/202401/c4bBhKzu1706685824726.jpeg'

$ws4.Range("B41").Value = '2024-05-10'
$ws4.Range("C41").Value = '广州·国际潮宠展—潮流创新宠物展会'
$ws4.Range("D41").Value = '阅江中路18号 广交会展馆C区'
$ws4.Range("E41").Value = '2024.05.10 10:30-05.12 18:30'
$ws4.Range("F41").Value = 54
$ws4.Range("G41").Value = 36
$ws4.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=82038'
$ws4.Range("I41").Value = '//i2.hdslb.com/bfs/openplatform/202402/om8irfxN1708678341525.jpeg'

$ws4.Range("B42").Value = '2024-05-18'
$ws4.Range("C42").Value = '广州·恋与深空only'
$ws4.Range("D42").Value = '大石街石北工业大道644号 巨大创意产业园'
$ws4.Range("E42").Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Range("F42").Value = 1598
$ws4.Range("G42").Value = 60
$ws4.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=81962'
$ws4.Range("I42").Value = '//i0.hdslb.com/bfs/openplatform/202402/a7aqaXrK1708485268977.jpeg'

$ws4.Range("B43").Value = '2024-05-18'
$ws4.Range("C43").Value = '广州·第五人格ONLY'
$ws4.Range("D43").Value = '洛浦街厦滘西环路1号 广州市岭南国际电子商务会展中心'
$ws4.Range("E43").Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Range("F43").Value = 329
$ws4.Range("G43").Value = 60
$ws4.Range("H43").Value = 'https://show.bilibili.com/platform/detail.html?id=82458'
$ws4.Range("I43").Value = '//i1.hdslb.com/bfs/openplatform/202403/D8jK0O2X1709778592031.jpeg'

$ws4.Range("B44").Value = '2024-05-25'
$ws4.Range("C44").Value = '广州·奶司的小人国娃展Nice Mini World  '
$ws4.Range("D44").Value = '洛浦街厦滘西环路1号 岭南会展中心'
$ws4.Range("E44").Value = '2024.05.25 10:30-05.25 17:00'
$ws4.Range("F44").Value = 29
$ws4.Range("G44").Value = 60
$ws4.Range("H44").Value = 'https://show.bilibili.com/platform/detail.html?id=82093'
$ws4.Range("I44").Value = '//i2.hdslb.com/bfs/openplatform/202402/rhIj7fnH1708936497981.jpeg'

# --- Sheet 4: drop the now-obsolete last row (was row 45) ---
$ws4.Rows.Item(45).Delete()